# Auto-generated Excel COM-interop edit script
# Applies the gh-pages data refresh (commit 456a3b4) to sheets 1 (展览) and 4 (全部类型)
$wb = $excel.ActiveWorkbook

# ---- Worksheet index 1 ----
$ws = $wb.Worksheets.Item(1)

# Bump view-count (F column) figures that changed for existing rows 1-13
$ws.Range('F3').Value = 566
$ws.Range('F4').Value = 548
$ws.Range('F7').Value = 39
$ws.Range('F10').Value = 4
$ws.Range('F12').Value = 4437
$ws.Range('F11').Value = 4641

# Insert a new row at 14: existing rows 14..15 shift down by one
$ws.Rows.Item(14).Insert()

# Restore the serial-number style/value for the newly-inserted A14 cell
$ws.Range('A14').Value = 13
$ws.Range('A14').Font.Bold = $true
$ws.Range('A14').HorizontalAlignment = -4108
$ws.Range('A14').VerticalAlignment = -4160
$ws.Range('A14').Borders.LineStyle = 1

# Write the new event's details into row 14 (合肥·星月动漫游戏展)
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = '2024-10-06'
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value = '合肥·星月动漫游戏展'
$ws.Range('D14').Value = '灵石路与皇藏峪路交叉口西南10米安徽百事兴电气有限公司院内2栋厂房2层 兄弟篮球俱乐部'
$ws.Range('E14').Value = '2024.10.06 10:00-10.06 17:00'
$ws.Range('F14').Value = 0
$ws.Range('G14').Value = 45
$ws.Range('H14').Value = 'https://show.bilibili.com/platform/detail.html?id=91958'
$ws.Range('I14').Value = '//i2.hdslb.com/bfs/openplatform/202409/mgB8U6bN1725361649767.jpeg'

# Bump the 想去人数 (interest count) for the two rows that shifted down
$ws.Range('F15').Value = 20  # was row 14, now row 15
$ws.Range('F16').Value = 153  # was row 15, now row 16

# Row-insert keeps each shifted cell's own stored value, so the serial-number
# column (A) needs re-numbering for every row that moved down one slot
$ws.Range('A15').Value = 14
$ws.Range('A16').Value = 15

# ---- Worksheet index 4 ----
$ws = $wb.Worksheets.Item(4)

# Bump view-count (F column) figures that changed for existing rows 1-13
$ws.Range('F3').Value = 566
$ws.Range('F4').Value = 548
$ws.Range('F7').Value = 39
$ws.Range('F10').Value = 4
$ws.Range('F12').Value = 4437
$ws.Range('F11').Value = 4641

# Insert a new row at 14: existing rows 14..18 shift down by one
$ws.Rows.Item(14).Insert()

# Restore the serial-number style/value for the newly-inserted A14 cell
$ws.Range('A14').Value = 13
$ws.Range('A14').Font.Bold = $true
$ws.Range('A14').HorizontalAlignment = -4108
$ws.Range('A14').VerticalAlignment = -4160
$ws.Range('A14').Borders.LineStyle = 1

# Write the new event's details into row 14 (合肥·星月动漫游戏展)
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = '2024-10-06'
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value = '合肥·星月动漫游戏展'
$ws.Range('D14').Value = '灵石路与皇藏峪路交叉口西南10米安徽百事兴电气有限公司院内2栋厂房2层 兄弟篮球俱乐部'
$ws.Range('E14').Value = '2024.10.06 10:00-10.06 17:00'
$ws.Range('F14').Value = 0
$ws.Range('G14').Value = 45
$ws.Range('H14').Value = 'https://show.bilibili.com/platform/detail.html?id=91958'
$ws.Range('I14').Value = '//i2.hdslb.com/bfs/openplatform/202409/mgB8U6bN1725361649767.jpeg'

# Bump the 想去人数 (interest count) for the two rows that shifted down
$ws.Range('F15').Value = 20  # was row 14, now row 15
$ws.Range('F16').Value = 153  # was row 15, now row 16

# Row-insert keeps each shifted cell's own stored value, so the serial-number
# column (A) needs re-numbering for every row that moved down one slot
$ws.Range('A15').Value = 14
$ws.Range('A16').Value = 15
$ws.Range('A17').Value = 16
$ws.Range('A18').Value = 17
$ws.Range('A19').Value = 18

Write-Output 'edit complete'
